$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - rows 2,3,4,5,6,7,10 in column F
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 311
$ws1.Range("F3").Value = 49
$ws1.Range("F4").Value = 475
$ws1.Range("F5").Value = 4560
$ws1.Range("F6").Value = 351
$ws1.Range("F7").Value = 624
$ws1.Range("F10").Value = 190

# Sheet "全部类型" (All types) - same data duplicated, but the row that was
# row 10 on "展览" is row 11 here (one extra row precedes it)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 311
$ws4.Range("F3").Value = 49
$ws4.Range("F4").Value = 475
$ws4.Range("F5").Value = 4560
$ws4.Range("F6").Value = 351
$ws4.Range("F7").Value = 624
$ws4.Range("F11").Value = 190
